$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.396201729774475
$ws.Range("B1").Value = 1.675781965255737
$ws.Range("C1").Value = 6.985039234161377
$ws.Range("D1").Value = 1.900041818618774
$ws.Range("E1").Value = 0.8592122793197632
